$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add Wins / Losses / Ties in AD1:AF1 ---
# Copy formatting (bold, centered, bordered) from the existing last header
# cell (AC1) so the new header cells match the rest of the header row.
$headerSrc = $ws.Range("AC1")
$headerCols = @("AD1", "AE1", "AF1")
$headerLabels = @("Wins", "Losses", "Ties")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $ws.Range($headerCols[$i])
    $headerSrc.Copy($cell)
    $cell.Value = $headerLabels[$i]
}

# --- Data rows: fill AD2:AF46 with the season record (Wins, Losses, Ties) ---
$wins = 91
$losses = 71
$ties = 0

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
